# Update the "想去人数" (want-to-go / interested count) values in column F
# across the four worksheets, per the upstream data refresh
# ("Update gh-pages to output generated at 456a3b4").
# Only numeric values change; no formatting/structure changes are required.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 370
$ws.Range("F4").Value = 144
$ws.Range("F5").Value = 1275
$ws.Range("F6").Value = 213
$ws.Range("F7").Value = 2444
$ws.Range("F8").Value = 865
$ws.Range("F9").Value = 18449
$ws.Range("F11").Value = 1857
$ws.Range("F14").Value = 315
$ws.Range("F15").Value = 590
$ws.Range("F23").Value = 75

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 26
$ws.Range("F8").Value = 113
$ws.Range("F9").Value = 109
$ws.Range("F14").Value = 66
$ws.Range("F19").Value = 10

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 5853
$ws.Range("F3").Value = 548

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 5853
$ws.Range("F4").Value = 548
$ws.Range("F6").Value = 370
$ws.Range("F8").Value = 144
$ws.Range("F10").Value = 1275
$ws.Range("F12").Value = 213
$ws.Range("F14").Value = 26
$ws.Range("F15").Value = 2444
$ws.Range("F16").Value = 865
$ws.Range("F17").Value = 18449
$ws.Range("F20").Value = 113
$ws.Range("F21").Value = 113
$ws.Range("F22").Value = 1857
$ws.Range("F24").Value = 109
$ws.Range("F26").Value = 315
$ws.Range("F27").Value = 590
$ws.Range("F36").Value = 66
$ws.Range("F44").Value = 10
$ws.Range("F48").Value = 75
